# Update Name of Algo
# Apply updated KNN-imputed values to the result data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.773
$ws.Range("D3").Value = -7.771000000000001

$ws.Range("A21").Value = -19.983

$ws.Range("A23").Value = -20.157

$ws.Range("D24").Value = -7.651999999999999

$ws.Range("A25").Value = -21.78900000000001

$ws.Range("C27").Value = -13.055

$ws.Range("C31").Value = -13.379

$ws.Range("C39").Value = -12.847

$ws.Range("C48").Value = -11.1

$ws.Range("C51").Value = -11.11

$ws.Range("C52").Value = -11.601

$ws.Range("A53").Value = -21.931

$ws.Range("C55").Value = -13.513

$ws.Range("C56").Value = -12.996

$ws.Range("A57").Value = -22.165
$ws.Range("C57").Value = -12.302
$ws.Range("D57").Value = -8.555000000000001

$ws.Range("A59").Value = -22.358

$ws.Range("D61").Value = -7.708

$ws.Range("A69").Value = -21.703

$ws.Range("D70").Value = -7.215000000000001

$ws.Range("C73").Value = -12.668

$ws.Range("A79").Value = -20.901

$ws.Range("A83").Value = -21.938

$ws.Range("D86").Value = -8.241

$ws.Range("C89").Value = -10.813

$ws.Range("C90").Value = -12.482

$ws.Range("A93").Value = -21.381

$ws.Range("D98").Value = -8.397

$ws.Range("D100").Value = -8.361999999999998

$ws.Range("D102").Value = -7.794000000000001
